$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45107
$ws.Range("M2").Value = 200

# Row 3
$ws.Range("D3").Value = 45043
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1083

# Row 4
$ws.Range("D4").Value = 45086
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20500
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 1139
